# Mock Test Issues reverted by Waqar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OpsTracker")

# Item 32 (row 26): "Online interview of Pronay Dhargave on 16 Nov at 3 PM" -> Done, Prannay selected
$ws.Range("D26").Value = "Done"
$ws.Range("E26").Value = "Prannay is selected."

# Rename "Debasish" -> "Debashish" wherever it appears as the Owned-by value
for ($r = 1; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Debasish") {
        $cell.Value = "Debashish"
    }
}

$ws.Range("C15").Select()
